# Commit message: "added strain names run 2641 to 4019"
# The generic placeholder strain labels "KN99A" / "KN99alpha" in column E (strain)
# are replaced with the actual strain names "TDY450" / "TDY451".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# KN99A -> TDY450
$ws.Range("E2").Value = "TDY450"
$ws.Range("E5").Value = "TDY450"

# KN99alpha -> TDY451
$ws.Range("E4").Value = "TDY451"
$ws.Range("E6").Value = "TDY451"
$ws.Range("E10").Value = "TDY451"

# Move/keep the active selection on E10, matching the saved selection state.
$ws.Range("E10").Select()
